$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BBW")

# Income statement: Earnings Before Interest And Taxes - oldest period (J) became unavailable
$ws.Range("J21").Value = "NA"

# Balance sheet: most-recent period (D) figures now available (were "NA")
$ws.Range("D41").Value = 21500
$ws.Range("D43").Value = 8300
$ws.Range("D44").Value = 58100
$ws.Range("D45").Value = 13300
$ws.Range("D46").Value = 101200
$ws.Range("D48").Value = 77700
$ws.Range("D49").Value = 900
$ws.Range("D52").Value = 7700
$ws.Range("D54").Value = 187500
$ws.Range("D57").Value = 19000
$ws.Range("D59").Value = 35800
$ws.Range("D60").Value = 54800
$ws.Range("D62").Value = 20600
$ws.Range("D66").Value = 75400
$ws.Range("D72").Value = 55900
$ws.Range("D76").Value = 112100

# Cash flow statement: oldest period (J) figures became unavailable
$ws.Range("J83").Value = "NA"
$ws.Range("J94").Value = "NA"
$ws.Range("J100").Value = "NA"
$ws.Range("J101").Value = "NA"
